# Update Name of Algo
# Apply updated result values produced by the algorithm run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value  = -12.834
$ws.Range("C3").Value  = -12.482
$ws.Range("C5").Value  = -12.18
$ws.Range("D7").Value  = -7.151999999999999
$ws.Range("B9").Value  = 6.417999999999999
$ws.Range("D9").Value  = -7.901999999999999
$ws.Range("C11").Value = -12.917
$ws.Range("C12").Value = -12.628
$ws.Range("B13").Value = 6.308000000000001
$ws.Range("B16").Value = 5.808
$ws.Range("B18").Value = 5.91
$ws.Range("B20").Value = 6.308000000000001
$ws.Range("C21").Value = -12.463
$ws.Range("D21").Value = -7.729000000000001
